# Append a new data row (row 6) to the IBB "Bag" sentiment data sheet,
# mirroring the layout/format of the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (date number format on column A, etc.) from the
# previous row so the new row's style matches the rest of the table.
$ws.Range("A5:N5").Copy()
$ws.Range("A6:N6").PasteSpecial()

$ws.Cells.Item(6, 1).Value = 42611.887604166666
$ws.Cells.Item(6, 2).Value = 4
$ws.Cells.Item(6, 3).Value = 52
$ws.Cells.Item(6, 4).Value = 45
$ws.Cells.Item(6, 5).Value = 100
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 6844
$ws.Cells.Item(6, 8).Value = 10709
$ws.Cells.Item(6, 9).Value = 1215
$ws.Cells.Item(6, 10).Value = 145
$ws.Cells.Item(6, 11).Value = 127
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = "Bag"
